$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Latitude" label/value pairs for Yanchep and Dawsville, placed to
# the right of the existing Note column (F:I), as part of the
# "added latitudes for areas" work-allocation update.
$ws.Range("F2").Value = "Latitude Yanchep"
$ws.Range("G2").Value = -31.547102426001601
$ws.Range("H2").Value = "Latitude Dawsville"
$ws.Range("I2").Value = -32.634355872262297

# Wrap the new label cells like the other note/label cells in the sheet,
# and give them a light box so they read as a little
# "Label | Value | Label | Value" strip: F2 gets both edges (it sits
# between two plain cells), H2 only needs a left edge since I2 (the
# rightmost cell) doesn't need a divider after it.
$ws.Range("F2").WrapText = $true
$ws.Range("F2").Borders.Item(7).LineStyle = 1
$ws.Range("F2").Borders.Item(10).LineStyle = 1

$ws.Range("H2").WrapText = $true
$ws.Range("H2").Borders.Item(7).LineStyle = 1

# Move the selection/view over to the newly added data.
$ws.Range("I3").Select()
